$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '43.701.31'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +1.27%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.243.00'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +0.59%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '321.52'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +2.08%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '101.33'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.41%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.578'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -1.32%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.554'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.24%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '37.26'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +0.60%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0829'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.37%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '7.66'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.37%  '
$ws.Range('E13').Value = '  -1.67%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.585.58'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.75%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.853'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.88%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '14.14'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -1.00%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.244.77'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.76%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '43.603.21'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +1.34%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '13.72'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -5.30%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0₃0983'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.88%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.43'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.95%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '65.16'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.71%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '3.16'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.09%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '235.99'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.60%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.15'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.68%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.05'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.96%  '
$ws.Range('B28').Value = 'InjectiveProtocol'
$ws.Range('C28').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '36.89'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +7.13%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.13'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -4.37%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '6.27'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -2.39%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '159.71'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +3.45%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '20.14'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -1.52%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0850'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -2.50%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.69'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -2.24%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.17'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +2.72%  '
$ws.Range('E36').Value = '  +8.78%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.92'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +0.68%  '
$ws.Range('E38').Value = '  -1.68%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.78'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +1.93%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '4.29'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -2.90%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '15.38'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +20.98%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0317'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -1.76%  '
$ws.Range('E43').Value = '  +0.29%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.809.21'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.84%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.200'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -2.38%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '82.66'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -5.09%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '5.24'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.43%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.71'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +5.83%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '74.30'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -4.08%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '58.72'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -1.07%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '103.60'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.31%  '
